$p = $ppt.ActivePresentation

# -----------------------------------------------------------------------
# 1) Slide 24 - "TextBox 9": fix typo "sfutta" -> "sfrutta"
# -----------------------------------------------------------------------
$slide24 = $p.Slides.Item(24)
$tbShape = $slide24.Shapes.Item(7)              # "TextBox 9"
$tbRange = $tbShape.TextFrame.TextRange
$fullText = $tbRange.Text
$pos = $fullText.IndexOf("sfutta")
if ($pos -ge 0) {
    $typo = $tbRange.Characters($pos + 1, 6)    # 1-based index, "sfutta" is 6 chars
    $typo.Text = "sfrutta"
}

# -----------------------------------------------------------------------
# 2) Slide 6 - "CasellaDiTesto 6": restyle "codice sorgente" and drop the
#    hyperlink, splitting the trailing run so ") " stays separate from
#    "durante le fasi di progettazione e codifica."
# -----------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$boxShape = $slide6.Shapes.Item(4)              # "CasellaDiTesto 6"
$boxRange = $boxShape.TextFrame.TextRange

$hlText = "codice sorgente"
$startPos = $boxRange.Text.IndexOf($hlText) + 1  # 1-based

# Remove the old hyperlinked run entirely (clears its text so the run,
# including its <a:hlinkClick/>, is dropped from the paragraph).
$hlRun = $boxRange.Characters($startPos, $hlText.Length)
$hlRun.Text = ""

# Grab an anchor made of the two characters that follow where the
# hyperlink used to be (") ") and insert the replacement text in front
# of it; the new text inherits the anchor run's plain formatting
# (solid bg1 fill, no hyperlink).
$anchor = $boxRange.Characters($startPos, 2)
$anchor.InsertBefore($hlText)

# Re-apply distinct formatting only to "codice sorgente" so it splits
# out into its own run: underline + Accent 6 (Lighter 40%).
$newRun = $boxRange.Characters($startPos, $hlText.Length)
$newRun.Font.Underline = $true
$newRun.Font.Color.ObjectThemeColor = 10        # msoThemeColorAccent6
$newRun.Font.Color.Brightness = 0.4             # "Lighter 40%" (lumMod 60000 / lumOff 40000)
